$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '96.647.89'
$ws.Range('E2').Value = '  +1.18%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.590.03'
$ws.Range('E3').Value = '  -0.25%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.97'
$ws.Range('E5').Value = '  +0.77%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '658.63'
$ws.Range('E6').Value = '  +0.98%  '

# Row 7
$ws.Range('E7').Value = '  +6.37%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.407'
$ws.Range('E8').Value = '  -0.49%  '

# Row 9
$ws.Range('E9').Value = '  +0.03%  '

# Row 10
$ws.Range('E10').Value = '  +4.20%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.589.94'
$ws.Range('E11').Value = '  -0.16%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.26'
$ws.Range('E12').Value = '  +0.26%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.202'
$ws.Range('E13').Value = '  +1.16%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.39'
$ws.Range('E14').Value = '  +1.24%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.254.05'
$ws.Range('E15').Value = '  -0.41%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.509.09'
$ws.Range('E16').Value = '  +1.28%  '

# Row 17
$ws.Range('E17').Value = '  +0.61%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.590.34'
$ws.Range('E18').Value = '  -0.29%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.74'
$ws.Range('E19').Value = '  -2.64%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.62'
$ws.Range('E20').Value = '  +0.88%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.84'
$ws.Range('E21').Value = '  -1.12%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.494'
$ws.Range('E22').Value = '  +1.88%  '

# Row 23
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.45'
$ws.Range('E23').Value = '  -1.16%  '

# Row 24
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '513.16'
$ws.Range('E24').Value = '  +0.31%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000202'
$ws.Range('E25').Value = '  +3.04%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.83'
$ws.Range('E26').Value = '  +2.19%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '96.77'
$ws.Range('E27').Value = '  +0.40%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.77'
$ws.Range('E28').Value = '  -1.23%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.782.15'
$ws.Range('E29').Value = '  -0.31%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.02'
$ws.Range('E30').Value = '  -3.87%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.149'
$ws.Range('E31').Value = '  +7.95%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.51'
$ws.Range('E32').Value = '  +1.82%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.01%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.187'
$ws.Range('E34').Value = '  +6.05%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.992'
$ws.Range('E35').Value = '  -0.28%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.71'
$ws.Range('E36').Value = '  -0.52%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.566'
$ws.Range('E37').Value = '  +1.25%  '

# Row 38
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '598.52'
$ws.Range('E38').Value = '  +6.61%  '

# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.50'
$ws.Range('E39').Value = '  +2.72%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.60'
$ws.Range('E40').Value = '  +9.07%  '

# Row 41
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.151'
$ws.Range('E41').Value = '  +0.62%  '

# Row 42
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.09%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.909'
$ws.Range('E43').Value = '  -1.89%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.84'
$ws.Range('E44').Value = '  +6.57%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.74'
$ws.Range('E45').Value = '  +0.75%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '34.47'
$ws.Range('E46').Value = '  +3.50%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.28'
$ws.Range('E47').Value = '  +0.86%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0421'
$ws.Range('E48').Value = '  +0.57%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.53'
$ws.Range('E49').Value = '  -0.96%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.62'
$ws.Range('E50').Value = '  +4.94%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.25'
$ws.Range('E51').Value = '  +1.98%  '
